# Deployed a2ecb74 to 0.5 with MkDocs 1.1.2 and mike 1.0.0
#
# 1) The cached "today" text of every date placeholder (slide master, all
#    11 slide layouts, notes master) moves from 5/12/21 -> 5/18/21.
# 2) The big white "Rectangle 4" placeholder shape on slide 1 is resized.
# 3) The footer / slide-number placeholder shapes are removed from slide 1.

$p = $ppt.ActivePresentation

function Update-DateText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "5/12/21") {
                $tr.Text = "5/18/21"
            }
        }
    }
}

# --- 1a) Slide master's own Date Placeholder ---
Update-DateText $p.SlideMaster.Shapes

# --- 1b) Every custom (slide) layout's Date Placeholder ---
$customLayouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $customLayouts.Count; $li++) {
    Update-DateText $customLayouts.Item($li).Shapes
}

# --- 1c) Notes master's Date Placeholder ---
Update-DateText $p.NotesMaster.Shapes

# --- 2) Resize/move "Rectangle 4" on slide 1 ---
$s = $p.Slides.Item(1)
$rect = $s.Shapes.Item(1)
$rect.Left = 101.34354400634766
$rect.Top = 63.72417449951172
$rect.Width = 703.0702514648438
$rect.Height = 322.88592529296875

# --- 3) Remove the footer & slide-number placeholders from slide 1 ---
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Footer Placeholder 1" -or $sh.Name -eq "Slide Number Placeholder 3") {
        $sh.Cut()
    }
}
